$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Work from the bottom of the table upward so that row indices used below
# always refer to their original (as-yet-unshifted) rows.

# Row 46: "10\t...\t100.0" -> "116"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "116"

# Row 45: "5\t...\t100.0" -> "0"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0"

# Row 44: "44\t...\t100.0" -> "100"
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"

# Row 12: "0.00011" -> "0.00006"
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.00006"

# Insert a new row after row 12 (before the old row 13) with "0.00247"
$newRow12 = $t.Rows.Add($t.Rows.Item(13))
$newRow12.Cells.Item(1).Range.Text = "0.00247"

# Row 11: "0.00004" -> "0.00005"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00005"

# Delete row 9 ("0.00003") entirely
$t.Rows.Item(9).Delete()

# Delete row 8 ("0.00000") entirely
$t.Rows.Item(8).Delete()

# Row 7: "0.00004" -> "0.00001"
$t.Rows.Item(7).Cells.Item(1).Range.Text = "0.00001"

# Insert a new row after row 5 (before the old row 6) with "0.00006"
$newRow5 = $t.Rows.Add($t.Rows.Item(6))
$newRow5.Cells.Item(1).Range.Text = "0.00006"

# Row 5: "0.00003" -> "0.00002"
$t.Rows.Item(5).Cells.Item(1).Range.Text = "0.00002"

# Row 4: "3" -> "62"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "62"

# Row 3: "116" -> "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# Row 2: "0" -> "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"

# Row 1: "100" -> "0M"
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
